$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Newmont Corporation)
$ws.Range("E2").Value = 61.8
$ws.Range("G2").Value = 60
$ws.Range("K2").Value = 73
$ws.Range("N2").Value = 85.87127175646313

# Row 3 (StreetTRACKS Gold Shares)
$ws.Range("E3").Value = 69.59999999999999
$ws.Range("G3").Value = 60
$ws.Range("K3").Value = 71.8
$ws.Range("N3").Value = 85.87127175646313

# Row 4 (Gold Dec 25)
$ws.Range("D4").Value = 4254.9
$ws.Range("E4").Value = 72.40000000000001
$ws.Range("F4").Value = 4.77
$ws.Range("I4").Value = 63
$ws.Range("J4").Value = 70
$ws.Range("K4").Value = 63
$ws.Range("N4").Value = 85.87127175646313
